$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 1.95
$ws.Range("M2").Value = 1.07
$ws.Range("N2").Value = 9
$ws.Range("AO2").Value = 7
$ws.Range("AU2").Value = 10
